$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The fixture previously padded the sheet with 5 empty, styled-only rows
# (rows 6-10). Remove them - the new test scenario only needs one extra
# data row (row 5) that duplicates the payment_id already used in row 4,
# to exercise the "duplicate payment id" validation.
$ws.Range("A6:K10").EntireRow.Delete()

# Build row 5 as a copy of row 4 (same formatting), then overwrite A5 so
# it duplicates row 4's payment_id on purpose.
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)

$ws.Range("H4:J4").Copy()
$ws.Range("H5:J5").PasteSpecial(-4122)

$ws.Range("H4").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").NumberFormat = $ws.Range("K4").NumberFormat

$ws.Range("A5").Value = $ws.Range("A4").Value2
$ws.Range("B5:G5").Value = $ws.Range("B4:G4").Value2
$ws.Range("H5:J5").Value = $ws.Range("H4:J4").Value2
$ws.Range("K5").Value = $ws.Range("K4").Value2
